$d = $word.ActiveDocument

# Replace every occurrence of "left-top" with "left-middle" (margin legend labels).
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute("left-top", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "left-middle", 2)

# Add a footer distance of 0.5" (36pt = 720 twips) to the section's page margins.
$d.PageSetup.FooterDistance = 36
